$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 879, shifting existing rows (879..957) down to (880..958)
$ws.Rows.Item(879).Insert()

# Populate the newly inserted row 879 with the new weekly data entry.
# Columns A,B,C,E,F,G,H,N,Q,R are constant across this data block, so copy them
# from the row directly below (old row 879, now at row 880).
$ws.Cells.Item(879, 1).Value2 = $ws.Cells.Item(880, 1).Value2   # A - Mercado ID
$ws.Cells.Item(879, 2).Value2 = $ws.Cells.Item(880, 2).Value2   # B - Mercado
$ws.Cells.Item(879, 3).Value2 = $ws.Cells.Item(880, 3).Value2   # C - Region
$ws.Cells.Item(879, 4).Value2 = 45132                            # D - Fecha
$ws.Cells.Item(879, 4).NumberFormat = $ws.Cells.Item(880, 4).NumberFormat
$ws.Cells.Item(879, 5).Value2 = $ws.Cells.Item(880, 5).Value2   # E - Codreg
$ws.Cells.Item(879, 6).Value2 = $ws.Cells.Item(880, 6).Value2   # F - Categoria ID
$ws.Cells.Item(879, 7).Value2 = $ws.Cells.Item(880, 7).Value2   # G - Categoria
$ws.Cells.Item(879, 8).Value2 = $ws.Cells.Item(880, 8).Value2   # H - Variedad
$ws.Cells.Item(879, 9).Value2 = "1a (guarda)"                   # I - Calidad
$ws.Cells.Item(879, 10).Value2 = 1600                            # J - Volumen
$ws.Cells.Item(879, 11).Value2 = 650                             # K - Precio minimo
$ws.Cells.Item(879, 12).Value2 = 700                             # L - Precio maximo
$ws.Cells.Item(879, 13).Value2 = 675                             # M - Precio promedio ponderado
$ws.Cells.Item(879, 14).Value2 = $ws.Cells.Item(880, 14).Value2  # N - Unidad de comercializacion
$ws.Cells.Item(879, 15).Value2 = "Región de O'Higgins"           # O - Origen
$ws.Cells.Item(879, 16).Value2 = 675                             # P - Precio $/Kg
$ws.Cells.Item(879, 17).Value2 = $ws.Cells.Item(880, 17).Value2  # Q - Kg o Unidades
$ws.Cells.Item(879, 18).Value2 = $ws.Cells.Item(880, 18).Value2  # R - Clasificacion
